# Add units to the formatted data headers and update the saved view.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New label text, introduced in the same order the original author's edit
# grew the shared-string table: dx2[cm] (in place of dx_2[cm]) first, then
# the appended unit-suffixed labels.
$ws.Range("C51").Value = "dx2[cm]"
$ws.Range("C19").Value = "Tges [s]"
$ws.Range("D19").Value = "T [s]"
$ws.Range("D44").Value = "T_s [s]"
$ws.Range("E19").Value = "f [Hz]"
$ws.Range("F19").Value = "w [rad/s]"
$ws.Range("D51").Value = "phi1 [°]"
$ws.Range("E51").Value = "phi2 [°]"

# Repeated block headers (rows 19, 26, 35, 58, 65, 74): n | Tges [s] | T [s] | f [Hz] | w [rad/s]
$headerRows = @(19, 26, 35, 58, 65, 74)
foreach ($r in $headerRows) {
    $ws.Range("C$r").Value = "Tges [s]"
    $ws.Range("D$r").Value = "T [s]"
    $ws.Range("E$r").Value = "f [Hz]"
    $ws.Range("F$r").Value = "w [rad/s]"
}

# Beat block headers (rows 44, 83): n | Tges [s] | T_s [s] | f [Hz] | w [rad/s]
$beatRows = @(44, 83)
foreach ($r in $beatRows) {
    $ws.Range("C$r").Value = "Tges [s]"
    $ws.Range("D$r").Value = "T_s [s]"
    $ws.Range("E$r").Value = "f [Hz]"
    $ws.Range("F$r").Value = "w [rad/s]"
}

# Coupling moment block headers (rows 51, 90): dx1[cm] | dx2[cm] | phi1 [°] | phi2 [°] | k
$couplingRows = @(51, 90)
foreach ($r in $couplingRows) {
    $ws.Range("C$r").Value = "dx2[cm]"
    $ws.Range("D$r").Value = "phi1 [°]"
    $ws.Range("E$r").Value = "phi2 [°]"
    $ws.Range("F$r").Value = "k"
}

# Saved scroll position moved up.
$ws.Application.ActiveWindow.ScrollRow = 62
